$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.083.56'
$ws.Range('E2').Value = '  -4.89%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.288.70'
$ws.Range('E3').Value = '  -5.39%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.24'
$ws.Range('E5').Value = '  -3.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.83'
$ws.Range('E6').Value = '  -3.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.594'
$ws.Range('E8').Value = '  -2.86%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.282.01'
$ws.Range('E9').Value = '  -5.28%  '
$ws.Range('E10').Value = '  -8.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.588'
$ws.Range('E11').Value = '  -4.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.62'
$ws.Range('E12').Value = '  -7.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000266'
$ws.Range('E14').Value = '  -5.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '635.54'
$ws.Range('E15').Value = '  -2.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.811.54'
$ws.Range('E16').Value = '  -5.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.129.91'
$ws.Range('E17').Value = '  -4.67%  '
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('E19').Value = '  -3.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.284.92'
$ws.Range('E20').Value = '  -5.45%  '
$ws.Range('E21').Value = '  -7.78%  '
$ws.Range('E22').Value = '  -4.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.25'
$ws.Range('E23').Value = '  +1.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '107.08'
$ws.Range('E24').Value = '  +8.09%  '
$ws.Range('E25').Value = '  -7.02%  '
$ws.Range('E26').Value = '  -7.38%  '
$ws.Range('E27').Value = '  -7.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.60'
$ws.Range('E28').Value = '  -3.57%  '
$ws.Range('E29').Value = '  -6.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.37'
$ws.Range('E30').Value = '  -6.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.01'
$ws.Range('E31').Value = '  -6.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.27'
$ws.Range('E32').Value = '  -6.69%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.05'
$ws.Range('E33').Value = '  -4.84%  '
$ws.Range('E34').Value = '  -3.52%  '
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '530.13'
$ws.Range('E35').Value = '  +0.47%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '57.60'
$ws.Range('E36').Value = '  -5.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.707.53'
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.998'
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('E39').Value = '  -4.37%  '
$ws.Range('D40').Value = '0.0₃0728'
$ws.Range('E40').Value = '  -7.88%  '
$ws.Range('B41').Value = 'CoreDAO'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.44'
$ws.Range('E41').Value = '  -1.94%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.73'
$ws.Range('E42').Value = '  -6.47%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.129'
$ws.Range('E43').Value = '  -3.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '33.07'
$ws.Range('E44').Value = '  -3.60%  '
$ws.Range('E45').Value = '  -9.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.28'
$ws.Range('E46').Value = '  -1.68%  '
$ws.Range('E47').Value = '  -5.81%  '
$ws.Range('E48').Value = '  -3.61%  '
$ws.Range('E49').Value = '  -7.64%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.998'
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('E51').Value = '  +3.13%  '
